$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Local / tiny.txt (45 Bytes) (10 times) ---
$ws.Range("B3").Value = "tiny.txt (45 Bytes) (10 times)"
$ws.Range("C3").Value = 0.003003400003944989
$ws.Range("D3").Value = 0.003274299998884089
$ws.Range("E3").Value = 0.006583900001714937
$ws.Range("F3").Value = -0.5026807822064099

# --- Row 4: Local / small.jpg (44.43 KB) (10 times) ---
$ws.Range("B4").Value = "small.jpg (44.43 KB) (10 times)"
$ws.Range("C4").Value = 0.04829529999988154
$ws.Range("D4").Value = 0.009827900004165713
$ws.Range("E4").Value = 0.04428600000101142
$ws.Range("F4").Value = -0.7780811090651388

# --- Row 5 (new): Local / medium.jpg (2.248 MB) (5 times) ---
$ws.Range("A5").Value = "Local"
$ws.Range("B5").Value = "medium.jpg (2.248 MB) (5 times)"
$ws.Range("C5").Value = 1.286583599995356
$ws.Range("D5").Value = 0.7791729999997188
$ws.Range("E5").Value = 0.5591157999995631
$ws.Range("F5").Value = 0.3935807215613074
$ws.Range("F5").HorizontalAlignment = -4131
$ws.Range("F5").IndentLevel = 1
$ws.Range("F5").NumberFormat = "[GREEN]0.00%;-[RED]0.00%"

# --- Row 6: now blank (previously held LAN / tiny.txt data) ---
$ws.Range("A6:I6").ClearContents()
$ws.Range("F6").ClearFormats()

# --- Row 7: Internet / tiny.txt (45 Bytes) (10 times) (previously LAN / medium.jpg) ---
$ws.Range("A7").Value = "Internet"
$ws.Range("B7").Value = "tiny.txt (45 Bytes) (10 times)"
$ws.Range("C7").Value = 0.2364815999972052
$ws.Range("D7").Value = 0.466005099999893
$ws.Range("E7").Value = 1.2279490999993867
$ws.Range("F7").Value = -0.6205012895077444

# --- Row 8 (new): Internet / small.jpg (44.43 KB) (10 times) ---
$ws.Range("A8").Value = "Internet"
$ws.Range("B8").Value = "small.jpg (44.43 KB) (10 times)"
$ws.Range("C8").Value = 0.3175075000021025
$ws.Range("D8").Value = 0.950518700000248
$ws.Range("E8").Value = 1.3214909000016632
$ws.Range("F8").Value = -0.2807224779231914
$ws.Range("F8").HorizontalAlignment = -4131
$ws.Range("F8").IndentLevel = 1
$ws.Range("F8").NumberFormat = "[GREEN]0.00%;-[RED]0.00%"

# --- Row 9: Internet / medium.jpg (2.248 MB) (5 times) (previously Internet / tiny.txt) ---
$ws.Range("B9").Value = "medium.jpg (2.248 MB) (5 times)"
$ws.Range("C9").Value = 3.2856530000048223
$ws.Range("D9").Value = 3.828223200002685
$ws.Range("E9").Value = 4.650853400002234
$ws.Range("F9").Value = -0.17687725869818938

# --- Row 10: now blank (previously held Internet / medium.jpg data) ---
$ws.Range("A10:I10").ClearContents()
$ws.Range("F10").ClearFormats()

Write-Host "edit complete"
